$d = $word.ActiveDocument

$d.Content.Find.Execute("641÷2=320, 1", $true, $false, $false, $false, $false, $true, 1, $false, "113÷7=16, 1", 2) | Out-Null
$d.Content.Find.Execute("800÷6=133, 2", $true, $false, $false, $false, $false, $true, 1, $false, "500÷4=125, 0", 2) | Out-Null
$d.Content.Find.Execute("945÷5=189, 0", $true, $false, $false, $false, $false, $true, 1, $false, "678÷3=226, 0", 2) | Out-Null
$d.Content.Find.Execute("381÷3=127, 0", $true, $false, $false, $false, $false, $true, 1, $false, "145÷6=24, 1", 2) | Out-Null
$d.Content.Find.Execute("873÷5=174, 3", $true, $false, $false, $false, $false, $true, 1, $false, "400÷2=200, 0", 2) | Out-Null
$d.Content.Find.Execute("359÷2=179, 1", $true, $false, $false, $false, $false, $true, 1, $false, "550÷4=137, 2", 2) | Out-Null
$d.Content.Find.Execute("900÷9=100, 0", $true, $false, $false, $false, $false, $true, 1, $false, "510÷5=102, 0", 2) | Out-Null
$d.Content.Find.Execute("755÷8=94, 3", $true, $false, $false, $false, $false, $true, 1, $false, "816÷2=408, 0", 2) | Out-Null
$d.Content.Find.Execute("893÷3=297, 2", $true, $false, $false, $false, $false, $true, 1, $false, "664÷4=166, 0", 2) | Out-Null
$d.Content.Find.Execute("783÷7=111, 6", $true, $false, $false, $false, $false, $true, 1, $false, "334÷9=37, 1", 2) | Out-Null
$d.Content.Find.Execute("757÷5=151, 2", $true, $false, $false, $false, $false, $true, 1, $false, "706÷5=141, 1", 2) | Out-Null
$d.Content.Find.Execute("878÷8=109, 6", $true, $false, $false, $false, $false, $true, 1, $false, "383÷9=42, 5", 2) | Out-Null
$d.Content.Find.Execute("422÷3=140, 2", $true, $false, $false, $false, $false, $true, 1, $false, "856÷4=214, 0", 2) | Out-Null
$d.Content.Find.Execute("979÷5=195, 4", $true, $false, $false, $false, $false, $true, 1, $false, "176÷6=29, 2", 2) | Out-Null
$d.Content.Find.Execute("712÷7=101, 5", $true, $false, $false, $false, $false, $true, 1, $false, "778÷9=86, 4", 2) | Out-Null
$d.Content.Find.Execute("316÷3=105, 1", $true, $false, $false, $false, $false, $true, 1, $false, "370÷9=41, 1", 2) | Out-Null
$d.Content.Find.Execute("387÷9=43, 0", $true, $false, $false, $false, $false, $true, 1, $false, "119÷8=14, 7", 2) | Out-Null
$d.Content.Find.Execute("456÷6=76, 0", $true, $false, $false, $false, $false, $true, 1, $false, "355÷8=44, 3", 2) | Out-Null
$d.Content.Find.Execute("684÷8=85, 4", $true, $false, $false, $false, $false, $true, 1, $false, "158÷8=19, 6", 2) | Out-Null
$d.Content.Find.Execute("109÷8=13, 5", $true, $false, $false, $false, $false, $true, 1, $false, "962÷2=481, 0", 2) | Out-Null
$d.Content.Find.Execute("394÷7=56, 2", $true, $false, $false, $false, $false, $true, 1, $false, "761÷7=108, 5", 2) | Out-Null
$d.Content.Find.Execute("406÷8=50, 6", $true, $false, $false, $false, $false, $true, 1, $false, "519÷3=173, 0", 2) | Out-Null
$d.Content.Find.Execute("557÷8=69, 5", $true, $false, $false, $false, $false, $true, 1, $false, "595÷3=198, 1", 2) | Out-Null
$d.Content.Find.Execute("745÷3=248, 1", $true, $false, $false, $false, $false, $true, 1, $false, "350÷2=175, 0", 2) | Out-Null
$d.Content.Find.Execute("828÷8=103, 4", $true, $false, $false, $false, $false, $true, 1, $false, "490÷8=61, 2", 2) | Out-Null
